# Update NATMI LR-pair stats (Hspg2-Ptprs) per Dr Hou's advice:
# ligand/receptor-expressing-cell counts move from 1 to 3, which
# changes every downstream average/total/specificity column (E,G:K,M:T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 135.7046583333333
$ws.Range("H2").Value = 407.113975
$ws.Range("I2").Value = 0.2901853119378819
$ws.Range("J2").Value = 0.2901853119378819
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.847811666666666
$ws.Range("N2").Value = 11.543435
$ws.Range("O2").Value = 0.0396810199351781
$ws.Range("P2").Value = 0.03968101993517809
$ws.Range("Q2").Value = 522.1659675560138
$ws.Range("R2").Value = 4699.493708004125
$ws.Range("S2").Value = 0.01151484914790297
$ws.Range("T2").Value = 0.01151484914790297

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 135.7046583333333
$ws.Range("H3").Value = 407.113975
$ws.Range("I3").Value = 0.2901853119378819
$ws.Range("J3").Value = 0.2901853119378819
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.80210733333333
$ws.Range("N3").Value = 173.406322
$ws.Range("O3").Value = 0.5960911739155557
$ws.Range("P3").Value = 0.5960911739155557
$ws.Range("Q3").Value = 7844.015226616661
$ws.Range("R3").Value = 70596.13703954994
$ws.Range("S3").Value = 0.1729769032461038
$ws.Range("T3").Value = 0.1729769032461038

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 135.7046583333333
$ws.Range("H4").Value = 407.113975
$ws.Range("I4").Value = 0.2901853119378819
$ws.Range("J4").Value = 0.2901853119378819
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.066157
$ws.Range("N4").Value = 18.198471
$ws.Range("O4").Value = 0.06255797260873913
$ws.Range("P4").Value = 0.06255797260873913
$ws.Range("Q4").Value = 823.2057630813583
$ws.Range("R4").Value = 7408.851867732225
$ws.Range("S4").Value = 0.01815340479566844
$ws.Range("T4").Value = 0.01815340479566844

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 135.7046583333333
$ws.Range("H5").Value = 407.113975
$ws.Range("I5").Value = 0.2901853119378819
$ws.Range("J5").Value = 0.2901853119378819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.25249166666667
$ws.Range("N5").Value = 87.757475
$ws.Range("O5").Value = 0.3016698335405271
$ws.Range("P5").Value = 0.301669833540527
$ws.Range("Q5").Value = 3969.699387023681
$ws.Range("R5").Value = 35727.29448321313
$ws.Range("S5").Value = 0.08754015474820677
$ws.Range("T5").Value = 0.08754015474820676

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 226.082006
$ws.Range("H6").Value = 678.246018
$ws.Range("I6").Value = 0.4834445496594812
$ws.Range("J6").Value = 0.4834445496594812
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.847811666666666
$ws.Range("N6").Value = 11.543435
$ws.Range("O6").Value = 0.0396810199351781
$ws.Range("P6").Value = 0.03968101993517809
$ws.Range("Q6").Value = 869.9209803102033
$ws.Range("R6").Value = 7829.28882279183
$ws.Range("S6").Value = 0.01918357281259107
$ws.Range("T6").Value = 0.01918357281259107

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 226.082006
$ws.Range("H7").Value = 678.246018
$ws.Range("I7").Value = 0.4834445496594812
$ws.Range("J7").Value = 0.4834445496594812
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.80210733333333
$ws.Range("N7").Value = 173.406322
$ws.Range("O7").Value = 0.5960911739155557
$ws.Range("P7").Value = 0.5960911739155557
$ws.Range("Q7").Value = 13068.01637694731
$ws.Range("R7").Value = 117612.1473925258
$ws.Range("S7").Value = 0.2881770291295973
$ws.Range("T7").Value = 0.2881770291295973

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 226.082006
$ws.Range("H8").Value = 678.246018
$ws.Range("I8").Value = 0.4834445496594812
$ws.Range("J8").Value = 0.4834445496594812
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.066157
$ws.Range("N8").Value = 18.198471
$ws.Range("O8").Value = 0.06255797260873913
$ws.Range("P8").Value = 0.06255797260873913
$ws.Range("Q8").Value = 1371.448943270942
$ws.Range("R8").Value = 12343.04048943848
$ws.Range("S8").Value = 0.03024331089544205
$ws.Range("T8").Value = 0.03024331089544205

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 226.082006
$ws.Range("H9").Value = 678.246018
$ws.Range("I9").Value = 0.4834445496594812
$ws.Range("J9").Value = 0.4834445496594812
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.25249166666667
$ws.Range("N9").Value = 87.757475
$ws.Range("O9").Value = 0.3016698335405271
$ws.Range("P9").Value = 0.301669833540527
$ws.Range("Q9").Value = 6613.461996498284
$ws.Range("R9").Value = 59521.15796848456
$ws.Range("S9").Value = 0.1458406368218508
$ws.Range("T9").Value = 0.1458406368218507

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1600446666666666
$ws.Range("H10").Value = 0.4801339999999999
$ws.Range("I10").Value = 0.0003422329939962955
$ws.Range("J10").Value = 0.0003422329939962955
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.847811666666666
$ws.Range("N10").Value = 11.543435
$ws.Range("O10").Value = 0.0396810199351781
$ws.Range("P10").Value = 0.03968101993517809
$ws.Range("Q10").Value = 0.6158217355877776
$ws.Range("R10").Value = 5.542395620289999
$ws.Range("S10").Value = 0.00001358015425724269
$ws.Range("T10").Value = 0.00001358015425724268

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1600446666666666
$ws.Range("H11").Value = 0.4801339999999999
$ws.Range("I11").Value = 0.0003422329939962955
$ws.Range("J11").Value = 0.0003422329939962955
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 57.80210733333333
$ws.Range("N11").Value = 173.406322
$ws.Range("O11").Value = 0.5960911739155557
$ws.Range("P11").Value = 0.5960911739155557
$ws.Range("Q11").Value = 9.250919000794221
$ws.Range("R11").Value = 83.25827100714798
$ws.Range("S11").Value = 0.0002040020671438871
$ws.Range("T11").Value = 0.0002040020671438871

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1600446666666666
$ws.Range("H12").Value = 0.4801339999999999
$ws.Range("I12").Value = 0.0003422329939962955
$ws.Range("J12").Value = 0.0003422329939962955
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.066157
$ws.Range("N12").Value = 18.198471
$ws.Range("O12").Value = 0.06255797260873913
$ws.Range("P12").Value = 0.06255797260873913
$ws.Range("Q12").Value = 0.9708560750126666
$ws.Range("R12").Value = 8.737704675113999
$ws.Range("S12").Value = 0.00002140940226422704
$ws.Range("T12").Value = 0.00002140940226422704

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1600446666666666
$ws.Range("H13").Value = 0.4801339999999999
$ws.Range("I13").Value = 0.0003422329939962955
$ws.Range("J13").Value = 0.0003422329939962955
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 29.25249166666667
$ws.Range("N13").Value = 87.757475
$ws.Range("O13").Value = 0.3016698335405271
$ws.Range("P13").Value = 0.301669833540527
$ws.Range("Q13").Value = 4.68170527796111
$ws.Range("R13").Value = 42.13534750165
$ws.Range("S13").Value = 0.0001032413703309387
$ws.Range("T13").Value = 0.0001032413703309386

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 105.7015583333333
$ws.Range("H14").Value = 317.104675
$ws.Range("I14").Value = 0.2260279054086406
$ws.Range("J14").Value = 0.2260279054086406
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.847811666666666
$ws.Range("N14").Value = 11.543435
$ws.Range("O14").Value = 0.0396810199351781
$ws.Range("P14").Value = 0.03968101993517809
$ws.Range("Q14").Value = 406.7196893398472
$ws.Range("R14").Value = 3660.477204058624
$ws.Range("S14").Value = 0.008969017820426818
$ws.Range("T14").Value = 0.008969017820426816

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 105.7015583333333
$ws.Range("H15").Value = 317.104675
$ws.Range("I15").Value = 0.2260279054086406
$ws.Range("J15").Value = 0.2260279054086406
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 57.80210733333333
$ws.Range("N15").Value = 173.406322
$ws.Range("O15").Value = 0.5960911739155557
$ws.Range("P15").Value = 0.5960911739155557
$ws.Range("Q15").Value = 6109.772820083927
$ws.Range("R15").Value = 54987.95538075535
$ws.Range("S15").Value = 0.1347332394727107
$ws.Range("T15").Value = 0.1347332394727107

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 105.7015583333333
$ws.Range("H16").Value = 317.104675
$ws.Range("I16").Value = 0.2260279054086406
$ws.Range("J16").Value = 0.2260279054086406
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.066157
$ws.Range("N16").Value = 18.198471
$ws.Range("O16").Value = 0.06255797260873913
$ws.Range("P16").Value = 0.06255797260873913
$ws.Range("Q16").Value = 641.2022479946584
$ws.Range("R16").Value = 5770.820231951925
$ws.Range("S16").Value = 0.01413984751536442
$ws.Range("T16").Value = 0.01413984751536442

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 105.7015583333333
$ws.Range("H17").Value = 317.104675
$ws.Range("I17").Value = 0.2260279054086406
$ws.Range("J17").Value = 0.2260279054086406
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 29.25249166666667
$ws.Range("N17").Value = 87.757475
$ws.Range("O17").Value = 0.3016698335405271
$ws.Range("P17").Value = 0.301669833540527
$ws.Range("Q17").Value = 3092.033954299514
$ws.Range("R17").Value = 27828.30558869562
$ws.Range("S17").Value = 0.06818580060013861
$ws.Range("T17").Value = 0.0681858006001386

